$wb = $excel.ActiveWorkbook

# --- Sheet "OFF" ---
$wsOff = $wb.Worksheets.Item("OFF")
$wsOff.Range("B2").Value = 280
$wsOff.Range("C2").Value = 189
$wsOff.Range("D2").Value = 54
$wsOff.Range("E2").Value = 28

# --- Sheet "DEF" ---
$wsDef = $wb.Worksheets.Item("DEF")
$wsDef.Range("B2").Value = 253
$wsDef.Range("C2").Value = 169
$wsDef.Range("D2").Value = 67
$wsDef.Range("E2").Value = 29
